# Switch license from BY-NC to BY-SA
# (units/8/lessons/4/resources/petascale-lesson-8.4-slides.pptx)
#
# The closing/license slide (slide 2) displays:
#   "CC BY-NC 4.0. To view a copy of this license, visit
#    https://creativecommons.org/licenses/by-nc/4.0"
# and should instead read:
#   "CC BY-SA 4.0. To view a copy of this license, visit
#    https://creativecommons.org/licenses/by-sa/4.0"
# with the hyperlink itself re-pointed at the by-sa license URL.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# 1. "CC BY-NC 4.0. ..." -> "CC BY-SA 4.0. ..."
$license = $tr.Find("BY-NC ", 0)
$license.Text = "BY-SA "

# 2. Displayed URL text: .../licenses/by-nc/4.0 -> .../licenses/by-sa/4.0
$urlText = $tr.Find("creativecommons.org/licenses/by-nc/4.0", 0)
$urlText.Text = "creativecommons.org/licenses/by-sa/4.0"

# 3. Underlying hyperlink target -> by-sa as well
$linkRun = $tr.Find("https://", 0)
$actionSettings = $linkRun.ActionSettings
$click = $actionSettings.Item(1)
$click.Hyperlink.Address = "https://creativecommons.org/licenses/by-sa/4.0"
